$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Curso MinTIC por tu mujer" title/casing, and the
# "... - Actualmente" dates to fixed end dates ("... - 2023")
$ws.Range("E12").Value = "Curso MinTIC Por Ti Mujer"
$ws.Range("B4").Value = "2017 - 2023"
$ws.Range("E4").Value = "Gestión de la comunicación (4 horas semanales - 2018 - 2023)"
$ws.Range("E5").Value = "Prácticas profesionales (4 horas semanales - 2022 - 2023)"

# Update selection to match the saved view state
$ws.Range("E5").Select()
